$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (also updates the _xlnm._FilterDatabase defined name's sheet ref)
$ws.Name = "PNEU"

# Update the 4 cells whose value changes from the (now removed) "primaire serie"
# shared string to the new quoted "Primaire serie" string.
$ws.Range("I10").Value = """Primaire serie"""
$ws.Range("I11").Value = """Primaire serie"""
$ws.Range("I15").Value = """Primaire serie"""
$ws.Range("I17").Value = """Primaire serie"""

# Remove the two data validations applied to column G.
$ws.Range("G28:G1048576").Validation.Delete()
$ws.Range("G2:G27").Validation.Delete()

# Shrink the AutoFilter to just the header row, extended to column I.
$ws.AutoFilterMode = $false
$ws.Range("A1:I1").AutoFilter(1)

# Keep the _xlnm._FilterDatabase defined name in sync with the new AutoFilter range.
$filterName = $wb.Names.Item("PNEU!_FilterDatabase")
$filterName.RefersTo = "=PNEU!`$A`$1:`$I`$1"

# Select the header row (A1:XFD1) instead of the previous H24 cell.
$ws.Range("A1:XFD1").Select()
